# Apply "Doing Updates for Financials" edit to HTBX_YR_FIN workbook.
# A new most-recent fiscal-year column of data was inserted at column D,
# shifting older data one column to the right (D->E->F->...); the furthest
# right column (J) either receives the old I value or becomes "NA" when the
# source row did not have enough historical data to shift into J.
# Values below are the exact final (post-edit) cell contents for every
# cell that changed, taken from the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("D8").Value = 5800
$ws.Range("E8").Value = 1500
$ws.Range("F8").Value = 300
$ws.Range("G8").Value = "NA"
$ws.Range("J8").Value = 0
# Row 12
$ws.Range("D12").Value = 16200
$ws.Range("E12").Value = 8300
$ws.Range("F12").Value = 9300
$ws.Range("G12").Value = 2600
$ws.Range("H12").Value = 2900
$ws.Range("I12").Value = 2700
$ws.Range("J12").Value = 900
# Row 14
$ws.Range("D14").Value = -500
$ws.Range("E14").Value = "NA"
$ws.Range("F14").Value = "NA"
$ws.Range("G14").Value = "NA"
$ws.Range("H14").Value = "NA"
$ws.Range("I14").Value = "NA"
$ws.Range("J14").Value = "NA"
# Row 17
$ws.Range("D17").Value = 23800
$ws.Range("E17").Value = 14900
$ws.Range("F17").Value = 13500
$ws.Range("G17").Value = 21000
$ws.Range("H17").Value = 12200
$ws.Range("I17").Value = 6600
$ws.Range("J17").Value = 2300
# Row 18
$ws.Range("D18").Value = -18000
$ws.Range("E18").Value = -13300
$ws.Range("F18").Value = -13100
$ws.Range("G18").Value = "NA"
$ws.Range("H18").Value = -12200
$ws.Range("I18").Value = -6600
$ws.Range("J18").Value = -2300
# Row 20
$ws.Range("D20").Value = 400
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 700
$ws.Range("G20").Value = "NA"
# Row 21
$ws.Range("D21").Value = -17400
$ws.Range("E21").Value = -13100
$ws.Range("F21").Value = -12300
$ws.Range("G21").Value = "NA"
$ws.Range("H21").Value = -12200
$ws.Range("I21").Value = -6500
$ws.Range("J21").Value = "NA"
# Row 22
$ws.Range("D22").Value = "NA"
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 500
$ws.Range("G22").Value = 400
# Row 23
$ws.Range("D23").Value = -17600
$ws.Range("E23").Value = -13200
$ws.Range("F23").Value = -13000
$ws.Range("G23").Value = -21100
$ws.Range("H23").Value = -12200
$ws.Range("I23").Value = -6600
$ws.Range("J23").Value = -2500
# Row 24
$ws.Range("D24").Value = -1000
$ws.Range("E24").Value = 7200
$ws.Range("G24").Value = "NA"
# Row 26
$ws.Range("D26").Value = "NA"
$ws.Range("E26").Value = -20400
$ws.Range("F26").Value = -13000
$ws.Range("G26").Value = -21100
$ws.Range("H26").Value = -12200
$ws.Range("I26").Value = -6600
$ws.Range("J26").Value = -2500
# Row 27
$ws.Range("D27").Value = "NA"
$ws.Range("E27").Value = -19900
$ws.Range("F27").Value = -12600
$ws.Range("G27").Value = -20300
$ws.Range("H27").Value = -11800
$ws.Range("I27").Value = -9100
$ws.Range("J27").Value = -2400
# Row 29
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = 8000
$ws.Range("H29").Value = "NA"
# Row 32
$ws.Range("D32").Value = -400
$ws.Range("E32").Value = -100
$ws.Range("F32").Value = -700
$ws.Range("G32").Value = "NA"
# Row 33
$ws.Range("D33").Value = "NA"
$ws.Range("E33").Value = -11800
$ws.Range("F33").Value = -12600
$ws.Range("G33").Value = -20300
$ws.Range("H33").Value = -11800
$ws.Range("I33").Value = -9100
$ws.Range("J33").Value = -2400
# Row 35
$ws.Range("D35").Value = "NA"
$ws.Range("E35").Value = -11800
$ws.Range("F35").Value = -12600
$ws.Range("G35").Value = -20300
$ws.Range("H35").Value = -11800
$ws.Range("I35").Value = -9100
$ws.Range("J35").Value = -2400
# Row 81
$ws.Range("D81").Value = "NA"
$ws.Range("E81").Value = -11800
$ws.Range("F81").Value = -12600
$ws.Range("G81").Value = -20300
$ws.Range("H81").Value = -11800
$ws.Range("I81").Value = -9100
$ws.Range("J81").Value = -2400
# Row 83
$ws.Range("J83").Value = "NA"
# Row 94
$ws.Range("J94").Value = "NA"
# Row 100
$ws.Range("J100").Value = "NA"
